{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfunction findParaStartingWith(prefix) {\n  for (const p of paras.items) {\n    if (p.text.indexOf(prefix) === 0) return p;\n  }\n  return null;\n}\n\nconst p1 = findParaStartingWith(\"Developed 17 of the 40 total field plug-ins\");\nconst p2 = findParaStartingWith(\"Developed and hosted 7 webinars\");\nconst p3 = findParaStartingWith(\"Wrote scripts and recorded audio\");\nconst p4 = findParaStartingWith(\"Wrote 57 user guides\");\nconst p5 = findParaStartingWith(\"Developed 23 Python scripts\");\n\nif (!p1 || !p2 || !p3 || !p4 || !p5) {\n  throw new Error(\"Could not find one of the target paragraphs: \" +\n    JSON.stringify([!!p1, !!p2, !!p3, !!p4, !!p5]));\n}\n\n// Capture the original (pre-edit) paragraph texts verbatim so whitespace\n// (e.g. the non-breaking space inside the webinars bullet) is preserved\n// exactly when that text is relocated to a new bullet position.\nconst origP1 = p1.text;\nconst origP2 = p2.text;\nconst origP3 = p3.text;\nconst origP4 = p4.text;\nconst origP5 = p5.text;\n\n// Paragraph 1: number/word edits + new trailing clause about encryption,\n// spliced into the original sentence.\nconst newP1 = origP1\n  .replace(\"Developed 17 of the 40 total field plug-ins\", \"Developed 18 of the 42 total public field plug-ins\")\n  .replace(\n    \"to enhance the capabilities of the software. Created custom\",\n    \"to enhance the capabilities of the software, as well as field plug-ins that integrate symmetric encryption. Created custom\"\n  );\n\n// Paragraph 2 (new) takes the old \"23 Python scripts\" bullet (originally last),\n// with the script count bumped and an encryption/decryption mention added.\nconst newP2 = origP5\n  .replace(\"Developed 23 Python scripts\", \"Developed 25 Python scripts\")\n  .replace(\"for data retrieval,\", \"for data encryption/decryption, data retrieval,\");\n\n// Paragraphs 3-5 (new) are simply the old paragraphs 2-4, unchanged, shifted\n// down by one bullet position.\nconst newP3 = origP2;\nconst newP4 = origP3;\nconst newP5 = origP4;\n\n// Apply replacements in document order (from the last paragraph backward is not required since\n// we operate on paragraph objects directly, not ranges, so order doesn't matter).\np1.insertText(newP1, \"Replace\");\np2.insertText(newP2, \"Replace\");\np3.insertText(newP3, \"Replace\");\np4.insertText(newP4, \"Replace\");\np5.insertText(newP5, \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the five target bullet paragraphs by their (unique) leading text.\n$p1 = $null; $p2 = $null; $p3 = $null; $p4 = $null; $p5 = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $t = $para.Range.Text\n    if ($t.StartsWith(\"Developed 17 of the 40 total field plug-ins\")) { $p1 = $para }\n    elseif ($t.StartsWith(\"Developed and hosted 7 webinars\")) { $p2 = $para }\n    elseif ($t.StartsWith(\"Wrote scripts and recorded audio\")) { $p3 = $para }\n    elseif ($t.StartsWith(\"Wrote 57 user guides\")) { $p4 = $para }\n    elseif ($t.StartsWith(\"Developed 23 Python scripts\")) { $p5 = $para }\n}\n\nif ($null -eq $p1 -or $null -eq $p2 -or $null -eq $p3 -or $null -eq $p4 -or $null -eq $p5) {\n    throw \"Could not find one of the target paragraphs\"\n}\n\n# Capture original paragraph text (paragraph mark included) verbatim before\n# any edits, so whitespace (e.g. non-breaking space) survives relocation.\n$origP1 = $p1.Range.Text\n$origP2 = $p2.Range.Text\n$origP3 = $p3.Range.Text\n$origP4 = $p4.Range.Text\n$origP5 = $p5.Range.Text\n\n# Strip the trailing paragraph-mark character Word includes in Range.Text.\n$mark = [char]13\nif ($origP1.EndsWith($mark)) { $origP1 = $origP1.Substring(0, $origP1.Length - 1) }\nif ($origP2.EndsWith($mark)) { $origP2 = $origP2.Substring(0, $origP2.Length - 1) }\nif ($origP3.EndsWith($mark)) { $origP3 = $origP3.Substring(0, $origP3.Length - 1) }\nif ($origP4.EndsWith($mark)) { $origP4 = $origP4.Substring(0, $origP4.Length - 1) }\nif ($origP5.EndsWith($mark)) { $origP5 = $origP5.Substring(0, $origP5.Length - 1) }\n\n# Paragraph 1: bump the counts and splice in the new clause about encryption.\n$newP1 = $origP1.Replace(\"Developed 17 of the 40 total field plug-ins\", \"Developed 18 of the 42 total public field plug-ins\")\n$newP1 = $newP1.Replace(\"to enhance the capabilities of the software. Created custom\", \"to enhance the capabilities of the software, as well as field plug-ins that integrate symmetric encryption. Created custom\")\n\n# Paragraph 2 (new): the old \"23 Python scripts\" bullet, moved up, with the\n# script count bumped and an encryption/decryption mention added.\n$newP2 = $origP5.Replace(\"Developed 23 Python scripts\", \"Developed 25 Python scripts\")\n$newP2 = $newP2.Replace(\"for data retrieval,\", \"for data encryption/decryption, data retrieval,\")\n\n# Paragraphs 3-5 (new) are simply the old paragraphs 2-4, unchanged, shifted\n# down by one bullet position.\n$newP3 = $origP2\n$newP4 = $origP3\n$newP5 = $origP4\n\n$p1.Range.Text = $newP1\n$p2.Range.Text = $newP2\n$p3.Range.Text = $newP3\n$p4.Range.Text = $newP4\n$p5.Range.Text = $newP5\n"}
